$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '37.445.50'
Set-TextValue $ws.Range("E2") '  +2.18%  '
Set-TextValue $ws.Range("D3") '2.037.86'
Set-TextValue $ws.Range("E3") '  +3.19%  '
Set-TextValue $ws.Range("E4") '  -0.09%  '
Set-TextValue $ws.Range("D5") '247.88'
Set-TextValue $ws.Range("E5") '  +1.56%  '
Set-TextValue $ws.Range("D6") '0.623'
Set-TextValue $ws.Range("E6") '  -0.64%  '
Set-TextValue $ws.Range("D7") '58.83'
Set-TextValue $ws.Range("E7") '  -2.15%  '
Set-TextValue $ws.Range("D9") '0.392'
Set-TextValue $ws.Range("E9") '  +3.62%  '
Set-TextValue $ws.Range("D10") '0.0809'
Set-TextValue $ws.Range("E10") '  +2.40%  '
Set-TextValue $ws.Range("E11") '  +0.22%  '
Set-TextValue $ws.Range("D12") '15.08'
Set-TextValue $ws.Range("E12") '  +5.76%  '
Set-TextValue $ws.Range("D13") '2.343.87'
Set-TextValue $ws.Range("E13") '  +3.56%  '
Set-TextValue $ws.Range("D14") '0.846'
Set-TextValue $ws.Range("E14") '  +0.05%  '
Set-TextValue $ws.Range("D15") '21.88'
Set-TextValue $ws.Range("E15") '  +1.05%  '
Set-TextValue $ws.Range("D16") '5.43'
Set-TextValue $ws.Range("E16") '  +2.44%  '
Set-TextValue $ws.Range("D17") '2.047.87'
Set-TextValue $ws.Range("E17") '  +4.10%  '
Set-TextValue $ws.Range("D18") '37.361.25'
Set-TextValue $ws.Range("E18") '  +2.04%  '
Set-TextValue $ws.Range("D19") '70.38'
Set-TextValue $ws.Range("E19") '  +0.96%  '
Set-TextValue $ws.Range("D20") '0.0₃0860'
Set-TextValue $ws.Range("E20") '  +0.52%  '
Set-TextValue $ws.Range("D21") '5.27'
Set-TextValue $ws.Range("E21") '  +3.36%  '
Set-TextValue $ws.Range("D22") '229.47'
Set-TextValue $ws.Range("E22") '  -0.13%  '
Set-TextValue $ws.Range("E23") '  -0.06%  '
Set-TextValue $ws.Range("E24") '  +4.24%  '
Set-TextValue $ws.Range("E25") '  -0.15%  '
Set-TextValue $ws.Range("D26") '9.30'
Set-TextValue $ws.Range("E26") '  +1.37%  '
Set-TextValue $ws.Range("D27") '164.26'
Set-TextValue $ws.Range("E27") '  +1.05%  '
Set-TextValue $ws.Range("D28") '0.137'
Set-TextValue $ws.Range("E28") '  -5.54%  '
Set-TextValue $ws.Range("D29") '19.90'
Set-TextValue $ws.Range("E29") '  +2.67%  '
Set-TextValue $ws.Range("D30") '1.35'
Set-TextValue $ws.Range("E30") '  +2.27%  '
Set-TextValue $ws.Range("D31") '0.122'
Set-TextValue $ws.Range("E31") '  +0.65%  '
Set-TextValue $ws.Range("D32") '0.0676'
Set-TextValue $ws.Range("E32") '  +10.01%  '
Set-TextValue $ws.Range("D33") '4.79'
Set-TextValue $ws.Range("E33") '  -0.88%  '
Set-TextValue $ws.Range("D34") '2.52'
Set-TextValue $ws.Range("E34") '  +10.51%  '
Set-TextValue $ws.Range("D35") '4.52'
Set-TextValue $ws.Range("E35") '  -0.60%  '
Set-TextValue $ws.Range("D36") '3.59'
Set-TextValue $ws.Range("E36") '  +5.94%  '
Set-TextValue $ws.Range("E37") '  -0.04%  '
Set-TextValue $ws.Range("E38") '  +2.35%  '
Set-TextValue $ws.Range("D39") '5.46'
Set-TextValue $ws.Range("E39") '  +0.50%  '
Set-TextValue $ws.Range("D40") '3.02'
Set-TextValue $ws.Range("E40") '  +3.58%  '
Set-TextValue $ws.Range("D41") '0.0981'
Set-TextValue $ws.Range("D42") '0.0219'
Set-TextValue $ws.Range("E42") '  +4.32%  '
Set-TextValue $ws.Range("E43") '  +1.04%  '
Set-TextValue $ws.Range("D44") '16.69'
Set-TextValue $ws.Range("E44") '  +4.65%  '
Set-TextValue $ws.Range("D45") '1.403.87'
Set-TextValue $ws.Range("E45") '  +2.44%  '
Set-TextValue $ws.Range("D46") '91.76'
Set-TextValue $ws.Range("E46") '  +2.98%  '
Set-TextValue $ws.Range("E47") '  +2.60%  '
Set-TextValue $ws.Range("D48") '7.49'
Set-TextValue $ws.Range("E48") '  +3.40%  '
Set-TextValue $ws.Range("D49") '2.12'
Set-TextValue $ws.Range("E49") '  +14.55%  '
Set-TextValue $ws.Range("D50") '2.88'
Set-TextValue $ws.Range("E50") '  +1.75%  '
Set-TextValue $ws.Range("D51") '2.235.20'
Set-TextValue $ws.Range("E51") '  +3.57%  '

Write-Host "Updated cryptos list"
